$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")

# Fix surcharge values - change from 1.0565 to 1 to match the other rows
$ws.Range("K16").Value = 1
$ws.Range("K17").Value = 1
$ws.Range("K20").Value = 1
$ws.Range("K21").Value = 1
$ws.Range("K25").Value = 1

# Update the active selection to C8 (support for longer quotes)
$ws.Activate()
$ws.Range("C8").Select()
